$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows 66 and 67 (columns P, AA, AB) ---
$ws.Range("P66").Value  = 796152
$ws.Range("AA66").Value = -79444
$ws.Range("AB66").Value = -1507792

$ws.Range("P67").Value  = 3522453
$ws.Range("AA67").Value = -77307
$ws.Range("AB67").Value = -4194398

# --- New row 68 (quarter 01-07-2021) ---
# Column A is a text "date-like" label ("01-07-2021"); force text so it is
# stored as a shared string instead of being auto-parsed into a date serial.
$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "01-07-2021"
$ws.Range("A68").Style = "Normal"

$rowValues = @{
    "B68" = 118912
    "C68" = 205751
    "D68" = 459691
    "E68" = 253940
    "F68" = -509018
    "G68" = 1502936
    "H68" = 2011954
    "I68" = 28689
    "J68" = 393490
    "K68" = 0
    "L68" = 0
    "M68" = 0
    "N68" = 0
    "O68" = 0
    "P68" = 8035464
    "Q68" = 5566866
    "R68" = 5928008
    "S68" = 5925987
    "T68" = 2021
    "U68" = 361142
    "V68" = 2555672
    "W68" = 4650816
    "X68" = 4650816
    "Y68" = 0
    "Z68" = 2095144
    "AA68" = -87075
    "AB68" = -7916552
    "AC68" = -3
}

foreach ($addr in $rowValues.Keys) {
    $ws.Range($addr).Value = $rowValues[$addr]
}

Write-Output "done"
